# Trade #357 closed at 2026-02-18 02:03:16 - unknown UNKNOWN +0.000%
#
# This script:
#   1) Updates the Summary sheet roll-up metrics.
#   2) Updates the Strategy Status row for MarketMaking.
#   3) Closes Trade #385 (MarketMaking, opened 01:52:19) in "All Trades" and
#      in the "MarketMaking" per-strategy sheet.
#   4) Appends four newly-opened trades (#414-#417) to "All Trades" and to
#      each of their respective per-strategy sheets
#      (momentum, HighProbConvergence, MarketMaking, EMAArbitrage).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1499.44   # Current Capital
$summary.Cells.Item(4, 2).Value = 0.55      # Total P&L $
$summary.Cells.Item(6, 2).Value = 385       # Total Trades
$summary.Cells.Item(7, 2).Value = 151       # Winning Trades
$summary.Cells.Item(9, 2).Value = 39.22     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(6, 3).Value = 99.42   # Capital
$status.Cells.Item(6, 4).Value = 238     # Trades
$status.Cells.Item(6, 5).Value = -0.39   # P&L $
$status.Cells.Item(6, 6).Value = -0.58   # P&L %
$status.Cells.Item(6, 7).Value = 36.13   # Win Rate %

# ---------------------------------------------------------------------
# 3) Close Trade #385 (MarketMaking) in "All Trades" (row 386)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(386, 7).Value = 0.87            # Exit Price
$allTrades.Cells.Item(386, 8).Value = "CLOSED"        # Status
$allTrades.Cells.Item(386, 9).Value = 2.3529          # P&L %
$allTrades.Cells.Item(386, 10).Value = 0.02           # P&L $
$allTrades.Cells.Item(386, 11).Value = 99.42          # Capital After
$allTrades.Cells.Item(386, 12).Value = "early_exit"   # Exit Reason
$allTrades.Cells.Item(386, 13).Value = 0.21           # Duration (min)

# ... and the matching row (239) in the "MarketMaking" per-strategy sheet
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(239, 7).Value = 0.87             # Exit Price
$marketMaking.Cells.Item(239, 8).Value = "CLOSED"         # Status
$marketMaking.Cells.Item(239, 9).Value = 2.3529           # P&L %
$marketMaking.Cells.Item(239, 10).Value = 0.02            # P&L $
$marketMaking.Cells.Item(239, 11).Value = 99.42           # Capital After
$marketMaking.Cells.Item(239, 16).Value = "early_exit"    # Exit Reason
$marketMaking.Cells.Item(239, 17).Value = 0.21            # Duration (min)

# ---------------------------------------------------------------------
# 4) Append newly opened trades
# ---------------------------------------------------------------------

# --- Trade #414 - momentum - DOWN --------------------------------------
$r = 415
$allTrades.Cells.Item($r, 1).Value = 414
$allTrades.Cells.Item($r, 2).Value = "'2026-02-18"
$allTrades.Cells.Item($r, 3).Value = "02:03:08"
$allTrades.Cells.Item($r, 4).Value = "momentum"
$allTrades.Cells.Item($r, 5).Value = "DOWN"
$allTrades.Cells.Item($r, 6).Value = 0.85
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 99.37699700270591
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.9
$allTrades.Cells.Item($r, 17).Value = "Downward momentum: -20.588% over 10 samples"

$momentum = $wb.Worksheets.Item("momentum")
$r = 81
$momentum.Cells.Item($r, 1).Value = 414
$momentum.Cells.Item($r, 2).Value = "'2026-02-18"
$momentum.Cells.Item($r, 3).Value = "02:03:08"
$momentum.Cells.Item($r, 4).Value = "momentum"
$momentum.Cells.Item($r, 5).Value = "DOWN"
$momentum.Cells.Item($r, 6).Value = 0.85
$momentum.Cells.Item($r, 8).Value = "OPEN"
$momentum.Cells.Item($r, 9).Value = 0
$momentum.Cells.Item($r, 10).Value = 0
$momentum.Cells.Item($r, 11).Value = 99.37699700270591
$momentum.Cells.Item($r, 12).Value = 0
$momentum.Cells.Item($r, 13).Value = 0
$momentum.Cells.Item($r, 14).Value = 0.9
$momentum.Cells.Item($r, 15).Value = "Downward momentum: -20.588% over 10 samples"
$momentum.Cells.Item($r, 17).Value = 0

# --- Trade #415 - HighProbConvergence - UP ------------------------------
$r = 416
$allTrades.Cells.Item($r, 1).Value = 415
$allTrades.Cells.Item($r, 2).Value = "'2026-02-18"
$allTrades.Cells.Item($r, 3).Value = "02:03:09"
$allTrades.Cells.Item($r, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item($r, 5).Value = "UP"
$allTrades.Cells.Item($r, 6).Value = 0.15
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 100.1931846556633
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.95
$allTrades.Cells.Item($r, 17).Value = "Mean reversion UP: price 19.76% below mean (z=-4.36)"

$highProb = $wb.Worksheets.Item("HighProbConvergence")
$r = 44
$highProb.Cells.Item($r, 1).Value = 415
$highProb.Cells.Item($r, 2).Value = "'2026-02-18"
$highProb.Cells.Item($r, 3).Value = "02:03:09"
$highProb.Cells.Item($r, 4).Value = "HighProbConvergence"
$highProb.Cells.Item($r, 5).Value = "UP"
$highProb.Cells.Item($r, 6).Value = 0.15
$highProb.Cells.Item($r, 8).Value = "OPEN"
$highProb.Cells.Item($r, 9).Value = 0
$highProb.Cells.Item($r, 10).Value = 0
$highProb.Cells.Item($r, 11).Value = 100.1931846556633
$highProb.Cells.Item($r, 12).Value = 0
$highProb.Cells.Item($r, 13).Value = 0
$highProb.Cells.Item($r, 14).Value = 0.95
$highProb.Cells.Item($r, 15).Value = "Mean reversion UP: price 19.76% below mean (z=-4.36)"
$highProb.Cells.Item($r, 17).Value = 0

# --- Trade #416 - MarketMaking - UP -------------------------------------
$r = 417
$allTrades.Cells.Item($r, 1).Value = 416
$allTrades.Cells.Item($r, 2).Value = "'2026-02-18"
$allTrades.Cells.Item($r, 3).Value = "02:03:10"
$allTrades.Cells.Item($r, 4).Value = "MarketMaking"
$allTrades.Cells.Item($r, 5).Value = "UP"
$allTrades.Cells.Item($r, 6).Value = 0.14
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 99.40221408909666
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.6
$allTrades.Cells.Item($r, 17).Value = "Normal spread capture: 247 bps"

$r = 247
$marketMaking.Cells.Item($r, 1).Value = 416
$marketMaking.Cells.Item($r, 2).Value = "'2026-02-18"
$marketMaking.Cells.Item($r, 3).Value = "02:03:10"
$marketMaking.Cells.Item($r, 4).Value = "MarketMaking"
$marketMaking.Cells.Item($r, 5).Value = "UP"
$marketMaking.Cells.Item($r, 6).Value = 0.14
$marketMaking.Cells.Item($r, 8).Value = "OPEN"
$marketMaking.Cells.Item($r, 9).Value = 0
$marketMaking.Cells.Item($r, 10).Value = 0
$marketMaking.Cells.Item($r, 11).Value = 99.40221408909666
$marketMaking.Cells.Item($r, 12).Value = 0
$marketMaking.Cells.Item($r, 13).Value = 0
$marketMaking.Cells.Item($r, 14).Value = 0.6
$marketMaking.Cells.Item($r, 15).Value = "Normal spread capture: 247 bps"
$marketMaking.Cells.Item($r, 17).Value = 0

# --- Trade #417 - EMAArbitrage - DOWN ------------------------------------
$r = 418
$allTrades.Cells.Item($r, 1).Value = 417
$allTrades.Cells.Item($r, 2).Value = "'2026-02-18"
$allTrades.Cells.Item($r, 3).Value = "02:03:10"
$allTrades.Cells.Item($r, 4).Value = "EMAArbitrage"
$allTrades.Cells.Item($r, 5).Value = "DOWN"
$allTrades.Cells.Item($r, 6).Value = 0.88
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 100.450616878256
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.7059
$allTrades.Cells.Item($r, 17).Value = "EMA:down, RSI:0.0, ROC:-20.59% | 2/3 DOWN"

$emaArb = $wb.Worksheets.Item("EMAArbitrage")
$r = 17
$emaArb.Cells.Item($r, 1).Value = 417
$emaArb.Cells.Item($r, 2).Value = "'2026-02-18"
$emaArb.Cells.Item($r, 3).Value = "02:03:10"
$emaArb.Cells.Item($r, 4).Value = "EMAArbitrage"
$emaArb.Cells.Item($r, 5).Value = "DOWN"
$emaArb.Cells.Item($r, 6).Value = 0.88
$emaArb.Cells.Item($r, 8).Value = "OPEN"
$emaArb.Cells.Item($r, 9).Value = 0
$emaArb.Cells.Item($r, 10).Value = 0
$emaArb.Cells.Item($r, 11).Value = 100.450616878256
$emaArb.Cells.Item($r, 12).Value = 0
$emaArb.Cells.Item($r, 13).Value = 0
$emaArb.Cells.Item($r, 14).Value = 0.7059
$emaArb.Cells.Item($r, 15).Value = "EMA:down, RSI:0.0, ROC:-20.59% | 2/3 DOWN"
$emaArb.Cells.Item($r, 17).Value = 0
